$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) — update column F (view/heat counters) for several rows
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 13889   # was 13886
$ws1.Range("F7").Value = 558     # was 557
$ws1.Range("F9").Value = 1785    # was 1784
$ws1.Range("F11").Value = 136    # was 135
$ws1.Range("F13").Value = 50     # was 49
$ws1.Range("F17").Value = 13946  # was 13939
$ws1.Range("F18").Value = 368    # was 367
$ws1.Range("F19").Value = 626    # was 625
$ws1.Range("F20").Value = 14968  # was 14967
$ws1.Range("F21").Value = 12     # was 11
$ws1.Range("F22").Value = 8267   # was 8264
$ws1.Range("F26").Value = 154    # was 153
$ws1.Range("F31").Value = 1039   # was 1038
$ws1.Range("F32").Value = 15     # was 13
$ws1.Range("F42").Value = 5089   # was 5088

# Sheet "全部类型" (sheet4) — same updates (mirrors the 展览 rows, minus F18)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 13889   # was 13886
$ws4.Range("F7").Value = 558     # was 557
$ws4.Range("F9").Value = 1785    # was 1784
$ws4.Range("F11").Value = 136    # was 135
$ws4.Range("F13").Value = 50     # was 49
$ws4.Range("F17").Value = 13946  # was 13939
$ws4.Range("F19").Value = 626    # was 625
$ws4.Range("F20").Value = 14968  # was 14967
$ws4.Range("F21").Value = 12     # was 11
$ws4.Range("F22").Value = 8267   # was 8264
$ws4.Range("F26").Value = 154    # was 153
$ws4.Range("F31").Value = 1039   # was 1038
$ws4.Range("F32").Value = 15     # was 14
$ws4.Range("F44").Value = 5089   # was 5088
